$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- F2: update rpc-reply message-id UUID ---
$f2 = $ws.Range("F2").Value()
$f2 = $f2 -replace "9f1c07fb-e75e-48a3-ae1a-b70201251b4d", "0b4d6dd8-568f-4368-8f65-7b6d15809b2b"
$ws.Range("F2").Value = $f2

# --- G2: update protocol identifier/name blocks ---
$g2 = $ws.Range("G2").Value()

$oldBlock1 = @'
<identifier>BGP</identifier>
              <name>BGP_65000</name>
'@

$newBlock1 = @'
<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
              <name>default</name>
'@

$oldBlock2 = @'
<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>BGP_65000</name>
'@

$newBlock2 = @'
<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>
                <name>default</name>
'@

$g2 = $g2.Replace($oldBlock1, $newBlock1)
$g2 = $g2.Replace($oldBlock2, $newBlock2)

$ws.Range("G2").Value = $g2
